# Added validations for Car on TripFolderPage
#
# The HotelScenarios data sheet gets a new "HotelName" value for the
# Registered-user / LAS Vegas basic-search scenario row (row 2), so the
# new hotel ("Holiday Inn Club Vacations LAS VEGAS") is available for the
# Car validations on TripFolderPage to reference.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GeneralScenarios")

# Fill in the previously-empty HotelName cell for row 2.
$ws.Range("J2").Value = "Holiday Inn Club Vacations LAS VEGAS"

# Leave the freeze-pane selection parked on the cell that was just edited,
# matching what Excel records after a user types into J2.
$ws.Range("J2").Select()
